$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows 1-2 shift down to 2-3
$ws.Rows.Item(1).Insert()

# New header row (row 1)
$ws.Range("A1").Value = "SORULAR"
$ws.Range("B1").Value = "GİTMESİ GEREKEN RAPOR"
$ws.Range("C1").Value = "GİTTİĞİ RAPOR"
$ws.Range("D1").Value = "EŞLEŞTİ Mİ?"
$ws.Range("E1").Value = "CEVAP"

# Row 2 (former row 1) - update the answer text
$ws.Range("E2").Value = "Bugün en yüksek POS girişi olan banka Ziraat Bankası, toplam girişi 590,040.30 TL."

# Row 3 (former row 2) - update matched report, match flag, and answer text
$ws.Range("C3").Value = "`n            `n                    📁 Günlük POS İşlemleri ve Banka Bazlı Toplamlar`n                  `n            `n          "
$ws.Range("D3").Value = "EVET"
$ws.Range("E3").Value = "Bugün bankalara göre en yüksek günlük giriş-çıkış farkı (net) AKBANK POS HS. ile -4,276,583.43 TL olarak kaydedilmiştir."

# The multi-line text in C3 can trigger an auto row-height bump; restore natural height
$ws.Rows.Item(3).AutoFit()

Write-Host $ws.UsedRange.Address()
